# Daily attendance processing - 2025-11-05 21:21:24
# Re-order the "Recorded By" (column G) value so that an originally-leading
# "System" entry is moved to the end of the comma-separated list, e.g.
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"
# Rows whose value does not start with "System, " (e.g. a lone "System",
# or a value that never had "System" first) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$prefix = "System, "

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -ne $null -and $text.StartsWith($prefix)) {
        $rest = $text.Substring($prefix.Length)
        $newValue = $rest + ", System"
        $cell.Value = $newValue
    }
}
